# Weekly data update: a new price-record row for "Espinaca" (Región
# Metropolitana, Mercado Mayorista Lo Valledor de Santiago) is inserted
# before the existing row 461, pushing the old rows 461-490 down to
# 462-491 (dimension grows from A1:R490 to A1:R491).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 461; this shifts rows 461..490 down
# to 462..491 and extends the used range automatically.
$ws.Rows.Item(461).Insert()

# Populate the new row 461 with the new weekly record.
$ws.Range("A461").Value2 = 6
$ws.Range("B461").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C461").Value2 = "Metropolitana"
$ws.Range("D461").Value2 = 44610
$ws.Range("E461").Value2 = 13
$ws.Range("F461").Value2 = 100112012
$ws.Range("G461").Value2 = "Espinaca"
$ws.Range("H461").Value2 = "Sin especificar"
$ws.Range("I461").Value2 = "Primera"
$ws.Range("J461").Value2 = 430
$ws.Range("K461").Value2 = 7500
$ws.Range("L461").Value2 = 8000
$ws.Range("M461").Value2 = 7698
$ws.Range("N461").Value2 = "$/cuna 10 kilos"
$ws.Range("O461").Value2 = "Región Metropolitana"
$ws.Range("P461").Value2 = 770
$ws.Range("Q461").Value2 = 10
$ws.Range("R461").Value2 = "Hortaliza"
